# Applies task-order stimulus filename corrections and sheet renames
# across the 5 worksheets of the workbook.

$wb = $excel.ActiveWorkbook

# --- Rename worksheets (new timestamped task-order identifiers) ---
$wb.Worksheets.Item(1).Name = "GNG_TO-1651255521342097"
$wb.Worksheets.Item(2).Name = "NB_TO-16512555235230985"
$wb.Worksheets.Item(3).Name = "RS_TO-16512555235311027"
$wb.Worksheets.Item(4).Name = "TOL_TO-16512555235870974"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16512555236660988"

# --- Sheet 1: GNG ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16512555213070977.csv"
$ws1.Range("B3").Value = "GNG_stims-16512555213250985.csv"
$ws1.Range("B4").Value = "go_stims-16512555213270981.csv"
$ws1.Range("B5").Value = "GNG_stims-16512555213410962.csv"

# --- Sheet 2: NB ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "ZB-match_0-1651255521362096.csv"
$ws2.Range("B3").Value = "ZB-match_8-16512555217950976.csv"
$ws2.Range("B4").Value = "TB-16512555234570985.csv"
$ws2.Range("B5").Value = "OB-1651255522778097.csv"
$ws2.Range("B6").Value = "TB-16512555233750966.csv"
$ws2.Range("B7").Value = "TB-16512555235110996.csv"
$ws2.Range("B8").Value = "ZB-match_1-16512555217391.csv"
$ws2.Range("B9").Value = "OB-16512555223520977.csv"
$ws2.Range("B10").Value = "OB-1651255522620096.csv"

# --- Sheet 3: RS ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = "eyes closed"
$ws3.Range("B3").Value = "eyes open"

# --- Sheet 4: TOL ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16512555235541015.csv"
$ws4.Range("B3").Value = "ZM_stims-16512555235331001.csv"
$ws4.Range("B4").Value = "MM_stims-16512555235700994.csv"
$ws4.Range("B5").Value = "ZM_stims-16512555235550988.csv"
$ws4.Range("B6").Value = "MM_stims-16512555235860965.csv"
$ws4.Range("B7").Value = "ZM_stims-16512555235711017.csv"

# --- Sheet 5: vSAT ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "vSAT_stims-16512555236340966.csv"
$ws5.Range("B3").Value = "vSAT_stims-16512555236500983.csv"
$ws5.Range("B4").Value = "SAT_stims-16512555235920978.csv"
$ws5.Range("B5").Value = "SAT_stims-16512555236180987.csv"
